# Updated symbol list on Fri Jan  6 09:23:13 UTC 2023 with GitHub Actions
#
# The "Price" (column D) and "Volume(1h)" (column E) figures on the crypto
# symbol sheet are refreshed with the latest scraped values. The cells hold
# plain text (not real numbers/percentages), so each value is written back
# as text (forcing a text number format while writing, then restoring the
# "Normal" style) to avoid Excel auto-converting strings like "256.84" or
# "-0.01%" into numeric/percentage values and losing formatting such as
# trailing zeros (e.g. "3.470", "0.0006050").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellAddress, $textValue)
    $range = $ws.Range($cellAddress)
    $range.NumberFormat = "@"
    $range.Value = $textValue
    $range.Style = "Normal"
}

# Row 2 (BNB)
Set-TextValue "D2" '256.84'
Set-TextValue "E2" '-0.01%'

# Row 3 (OKB)
Set-TextValue "D3" '26.88'
Set-TextValue "E3" '-0.90%'

# Row 4 (HuobiToken)
Set-TextValue "D4" '4.638'
Set-TextValue "E4" '-2.44%'

# Row 5 (Cronos)
Set-TextValue "D5" '0.05895'
Set-TextValue "E5" '-0.83%'

# Row 6 (KuCoinToken)
Set-TextValue "D6" '6.609'
Set-TextValue "E6" '-0.68%'

# Row 7 (MXToken) - price unchanged
Set-TextValue "E7" '-2.26%'

# Row 8 (FTXToken)
Set-TextValue "D8" '0.9201'
Set-TextValue "E8" '-3.16%'

# Row 9 (WazirX)
Set-TextValue "D9" '0.1377'
Set-TextValue "E9" '-1.70%'

# Row 10 (LiechtensteinCryptoassetsExchange)
Set-TextValue "D10" '0.04286'
Set-TextValue "E10" '8.28%'

# Row 11 (MandalaExchangeToken)
Set-TextValue "D11" '0.07002'
Set-TextValue "E11" '-2.16%'

# Row 12 (BitrueCoin)
Set-TextValue "D12" '0.03057'
Set-TextValue "E12" '-4.41%'

# Row 13 (BitMartToken)
Set-TextValue "D13" '0.09104'
Set-TextValue "E13" '-1.46%'

# Row 14 (BitForexToken)
Set-TextValue "D14" '0.001526'
Set-TextValue "E14" '-1.08%'

# Row 15 (One)
Set-TextValue "D15" '0.0006050'
Set-TextValue "E15" '-94.23%'

# Row 16 (TigerCash)
Set-TextValue "D16" '0.006013'
Set-TextValue "E16" '-0.25%'

# Row 17 (LEO)
Set-TextValue "D17" '3.470'
Set-TextValue "E17" '-0.37%'

# Row 18 (GateToken)
Set-TextValue "D18" '3.165'
Set-TextValue "E18" '-1.19%'

# Row 19 (BTSEToken) - price unchanged
Set-TextValue "E19" '-1.12%'

# Row 20 (BitpandaEcosystemToken) - price unchanged
Set-TextValue "E20" '-2.48%'

# Row 21 (ProBitToken) - price unchanged
Set-TextValue "E21" '-0.24%'

# Row 22 (MCDex)
Set-TextValue "D22" '3.904'
Set-TextValue "E22" '2.41%'

# Row 23 (CoinExToken)
Set-TextValue "D23" '0.04257'
Set-TextValue "E23" '1.11%'

# Row 24 (BitKan)
Set-TextValue "D24" '0.001224'
Set-TextValue "E24" '0.26%'

# Row 25 (HotbitToken)
Set-TextValue "D25" '0.004291'
Set-TextValue "E25" '-4.59%'

# Row 26 (NitroEx) - price unchanged
Set-TextValue "E26" '0.04%'

# Row 27 (UpBots) - price unchanged
Set-TextValue "E27" '-21.34%'

# Row 40 (IDEX)
Set-TextValue "D40" '0.03781'
Set-TextValue "E40" '-1.12%'

# Row 41 (KickToken)
Set-TextValue "D41" '0.006313'
Set-TextValue "E41" '1.82%'

# Row 42 (BKEXToken)
Set-TextValue "D42" '0.1098'
Set-TextValue "E42" '-0.24%'

# Row 43 (CEJI)
Set-TextValue "D43" '0.002201'
Set-TextValue "E43" '-2.22%'

# Row 44 (LocalTraders)
Set-TextValue "D44" '0.01412'
Set-TextValue "E44" '33.66%'

# Row 45 (CoinLion)
Set-TextValue "D45" '0.00005361'
Set-TextValue "E45" '-2.52%'

# Row 46 (Kangarootoken) - price unchanged
Set-TextValue "E46" '0.05%'

# Row 47 (CoinbaseStockToken)
Set-TextValue "D47" '0.04550'
Set-TextValue "E47" '-48.59%'

# Row 48 (BOLO) - price unchanged
Set-TextValue "E48" '10,478.06%'

# Row 49 (CryptobidCoin)
Set-TextValue "D49" '0.00002101'
Set-TextValue "E49" '0.05%'

# Row 50 (SpecialPowerGold)
Set-TextValue "D50" '0.0002001'
Set-TextValue "E50" '0.05%'
